# Sync attendance_reports, modules_schedules, and assets from main repo
#
# The "Recorded By" column (G) on the "Session Analysis Results" sheet
# lists the people/systems who recorded each attendance entry as a
# comma-separated string (e.g. "dnasr281@gmail.com, System"). Upstream
# re-ordered that list for a specific set of rows, swapping the two
# names around the comma. This script reproduces that swap for exactly
# the affected rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows in column G whose two comma-separated "Recorded By" values need
# to be swapped (order reversed), as identified from the commit diff.
$rows = @(
    3, 6, 7, 10, 11, 12, 13, 14, 15, 17, 18, 19, 20, 21, 22,
    30, 33, 34, 37, 38, 39, 40, 41, 42, 44, 45, 46, 47, 48, 49,
    57, 60, 61, 64, 65, 66, 67, 68, 69, 71, 72, 73, 74, 75, 76,
    86, 87, 88, 89, 90, 93, 95, 96, 97, 99, 102,
    112, 113, 114, 115, 116, 119, 121, 122, 123, 125, 128,
    138, 139, 140, 141, 142, 145, 147, 148, 149, 151, 154
)

foreach ($row in $rows) {
    $cell = $ws.Range("G$row")
    $current = [string]$cell.Value2
    $parts = $current -split ", "
    if ($parts.Length -eq 2) {
        $cell.Value = "$($parts[1]), $($parts[0])"
    }
}
